$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 3
$ws.Range("D3").Value = 10.34
$ws.Range("E3").Value = 10.72

# Row 4
$ws.Range("C4").Value = 9.66
$ws.Range("E4").Value = 10.68
$ws.Range("F4").Value = 10.13
$ws.Range("H4").Value = 8.67

# Row 5
$ws.Range("C5").Value = 9.23
$ws.Range("D5").Value = 9.32
$ws.Range("F5").Value = 10.2

# Row 6
$ws.Range("D6").Value = 9.87
$ws.Range("E6").Value = 9.8
$ws.Range("G6").Value = 10.3

# Row 7
$ws.Range("F7").Value = 9.7
$ws.Range("I7").Value = 7.08

# Row 8
$ws.Range("D8").Value = 11.33
$ws.Range("I8").Value = 8.86

# Row 9
$ws.Range("G9").Value = 12.92
$ws.Range("H9").Value = 11.14
